# Update "Pais" (countries) sheet with the latest COVID-19 snapshot.
#
# The underlying sharedStrings table in the source workbook is reordered
# upstream (new countries spliced into the ranked list, pushing the rows
# below them down by one), and a batch of case/death counters are refreshed
# for the affected rows. Rather than touch the shared-string table directly
# (an internal storage detail Excel/COM manages on its own), we simply set
# every cell to its final value; Excel re-points/creates shared-string
# entries as needed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Footer timestamp (row 1)
$ws.Range("A1").Value = "Datos actualizados a 23 de Marzo de 2020 a las 21:46"

$ws.Range("A1").Value = "Datos actualizados a 23 de Marzo de 2020 a las 21:46"

# Row 6: Estados Unidos - refreshed totals
$ws.Range("B6").Value = 42379
$ws.Range("C6").Value = 8833
$ws.Range("E6").Value = 41567
$ws.Range("G6").Value = 98
$ws.Range("H6").Value = 517

# Row 22: Brasil - refreshed totals; Turquia spliced in right below it
$ws.Range("B22").Value = 1696
$ws.Range("C22").Value = 150
$ws.Range("E22").Value = 1660
$ws.Range("G22").Value = 9
$ws.Range("H22").Value = 34

# Row 23: Turquia (new position) with its own totals
$ws.Range("A23").Value = "Turquia"
$ws.Range("B23").Value = 1529
$ws.Range("C23").Value = 293
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 1492
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 7
$ws.Range("H23").Value = 37

# Row 24: Malasia shifted down one row
$ws.Range("A24").Value = "Malasia"
$ws.Range("B24").Value = 1518
$ws.Range("C24").Value = 212
$ws.Range("D24").Value = 159
$ws.Range("E24").Value = 1345
$ws.Range("F24").Value = 57
$ws.Range("G24").Value = 4
$ws.Range("H24").Value = 14

# Row 25: Dinamarca shifted down one row
$ws.Range("A25").Value = "Dinamarca"
$ws.Range("B25").Value = 1450
$ws.Range("C25").Value = 55
$ws.Range("D25").Value = 1
$ws.Range("E25").Value = 1425
$ws.Range("F25").Value = 55
$ws.Range("G25").Value = 11
$ws.Range("H25").Value = 24

# Row 26: Israel shifted down one row
$ws.Range("A26").Value = "Israel"
$ws.Range("B26").Value = 1442
$ws.Range("C26").Value = 371
$ws.Range("D26").Value = 41
$ws.Range("E26").Value = 1400
$ws.Range("F26").Value = 29

# Row 27: Chequia shifted down one row
$ws.Range("A27").Value = "Chequia"
$ws.Range("C27").Value = 116
$ws.Range("D27").Value = 7
$ws.Range("E27").Value = 1228
$ws.Range("F27").Value = 19
$ws.Range("H27").Value = 1

# Row 64: Argelia..Taiwan block shifted down one row; Armenia spliced in
$ws.Range("A64").Value = "Armenia"
$ws.Range("B64").Value = 235
$ws.Range("C64").Value = 41
$ws.Range("D64").Value = 2
$ws.Range("E64").Value = 233
$ws.Range("F64").Value = 6
$ws.Range("H64").Value = 0

$ws.Range("A65").Value = "Argelia"
$ws.Range("B65").Value = 230
$ws.Range("C65").Value = 29
$ws.Range("D65").Value = 65
$ws.Range("E65").Value = 148
$ws.Range("F65").Value = 0
$ws.Range("H65").Value = 17

$ws.Range("A66").Value = "Bulgaria"
$ws.Range("B66").Value = 201
$ws.Range("C66").Value = 14
$ws.Range("D66").Value = 3
$ws.Range("E66").Value = 195
$ws.Range("F66").Value = 8
$ws.Range("H66").Value = 3

$ws.Range("A67").Value = "Emiratos Arabes Unidos"
$ws.Range("B67").Value = 198
$ws.Range("C67").Value = 45
$ws.Range("D67").Value = 41
$ws.Range("E67").Value = 155
$ws.Range("F67").Value = 2

$ws.Range("A68").Value = "Taiwan"
$ws.Range("B68").Value = 195
$ws.Range("C68").Value = 26
$ws.Range("D68").Value = 28
$ws.Range("E68").Value = 165
$ws.Range("F68").Value = 0
$ws.Range("H68").Value = 2

# Row 107: Trinidad yTobago / Liechtenstein swap back to original relative order
$ws.Range("A107").Value = "Trinidad yTobago"
$ws.Range("C107").Value = 1

$ws.Range("A108").Value = "Liechtenstein"
$ws.Range("C108").Value = 14

# Row 127: Guatemala spliced in above Guyana
$ws.Range("A127").Value = "Guatemala"

$ws.Range("A128").Value = "Guyana"

# Row 132: Togo spliced in above Polinesia Francesa
$ws.Range("A132").Value = "Togo"
$ws.Range("C132").Value = 2

$ws.Range("A133").Value = "Polinesia Francesa"
$ws.Range("C133").Value = 0

# Row 134: Islas Virgenes de los Estados Unidos moves above Barbados
$ws.Range("A134").Value = "Islas Virgenes de los Estados Unidos"
$ws.Range("C134").Value = 11

$ws.Range("A135").Value = "Barbados"
$ws.Range("C135").Value = 3

# Row 147: San Martin (Parte Francesa) moves above Nueva Caledonia
$ws.Range("A147").Value = "San Martin (Parte Francesa)"
$ws.Range("C147").Value = 3

$ws.Range("A148").Value = "Nueva Caledonia"
$ws.Range("C148").Value = 4

# Row 154: Congo/Suazilandia/Guinea spliced in above Groenlandia
$ws.Range("A154").Value = "Congo"
$ws.Range("C154").Value = 1

$ws.Range("A155").Value = "Suazilandia"
$ws.Range("C155").Value = 0

$ws.Range("A156").Value = "Guinea"
$ws.Range("C156").Value = 2

$ws.Range("A157").Value = "Groenlandia"
$ws.Range("C157").Value = 0

$ws.Range("A158").Value = "Namibia"
$ws.Range("C158").Value = 1

$ws.Range("A159").Value = "Bahamas"
$ws.Range("C159").Value = 0

# Row 161: Large African/Caribbean block reshuffled
$ws.Range("A161").Value = "Republica de Africa Central"

$ws.Range("A162").Value = "El Salvador"

$ws.Range("A163").Value = "Santa Lucia"

$ws.Range("A164").Value = "Fiyi"
$ws.Range("C164").Value = 1

$ws.Range("A165").Value = "Antigua y Barbuda"
$ws.Range("C165").Value = 2

$ws.Range("A166").Value = "Liberia"
$ws.Range("C166").Value = 0

$ws.Range("A167").Value = "Angola"
$ws.Range("C167").Value = 1

$ws.Range("A168").Value = "San Bartolome"
$ws.Range("C168").Value = 0

$ws.Range("A169").Value = "Zambia"
$ws.Range("C169").Value = 0

$ws.Range("A171").Value = "Republica de Yibuti"
$ws.Range("C171").Value = 2

# Row 174: Nicaragua..Benin block reshuffled ahead of San Martin (Parte Holandesa)
$ws.Range("A174").Value = "Nicaragua"
$ws.Range("C174").Value = 0

$ws.Range("A175").Value = "Birmania"
$ws.Range("C175").Value = 2

$ws.Range("A176").Value = "Butan"

$ws.Range("A177").Value = "Mauritania"
$ws.Range("C177").Value = 0

$ws.Range("A178").Value = "Niger"

$ws.Range("A179").Value = "Benin"

$ws.Range("A180").Value = "San Martin (Parte Holandesa)"
$ws.Range("C180").Value = 1

# Row 181: Nepal takes this row; small counter tweaks ripple through
$ws.Range("A181").Value = "Nepal"
$ws.Range("C181").Value = 1
$ws.Range("D181").Value = 1
$ws.Range("H181").Value = 0

$ws.Range("A182").Value = "Sudan"
$ws.Range("C182").Value = 0
$ws.Range("G182").Value = 0

$ws.Range("A183").Value = "Gambia"
$ws.Range("D183").Value = 0
$ws.Range("G183").Value = 1
$ws.Range("H183").Value = 1

# Row 184: Siria..Eritrea block reshuffled
$ws.Range("A184").Value = "Siria"
$ws.Range("C184").Value = 0

$ws.Range("A185").Value = "Granada"

$ws.Range("A186").Value = "Dominica"

$ws.Range("A187").Value = "Republica del Chad"

$ws.Range("A188").Value = "Mozambique"

$ws.Range("A189").Value = "Uganda"

$ws.Range("A190").Value = "Montserrat"

$ws.Range("A191").Value = "Eritrea"

# Row 192: Belice..San Vicente y las Granadinas tail reshuffled
$ws.Range("A192").Value = "Belice"
$ws.Range("C192").Value = 1

$ws.Range("A193").Value = "Islas Turcas y Caicos"
$ws.Range("C193").Value = 1

$ws.Range("A194").Value = "Timor Oriental"

$ws.Range("A195").Value = "Papua Nueva Guinea"

$ws.Range("A196").Value = "Somalia"
$ws.Range("C196").Value = 0

$ws.Range("A197").Value = "Santa Sede"

$ws.Range("A198").Value = "San Vicente y las Granadinas"
